$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue "D2" "26.291.08"
Set-TextValue "E2" "  -0.02%  "
Set-TextValue "D3" "1.689.66"
Set-TextValue "E3" "  +0.68%  "
Set-TextValue "D4" "1.008"
Set-TextValue "E4" "  -0.09%  "
Set-TextValue "D5" "217.74"
Set-TextValue "E5" "  -0.10%  "
Set-TextValue "D6" "0.5325"
Set-TextValue "E6" "  +1.26%  "
Set-TextValue "D7" "1.008"
Set-TextValue "E7" "  -0.10%  "
Set-TextValue "D8" "0.2719"
Set-TextValue "E8" "  +1.30%  "
Set-TextValue "D9" "0.06415"
Set-TextValue "E9" "  -0.75%  "
Set-TextValue "D10" "21.70"
Set-TextValue "E10" "  -0.80%  "
Set-TextValue "D11" "0.07691"
Set-TextValue "E11" "  +2.36%  "
Set-TextValue "D12" "1.703.95"
Set-TextValue "E12" "  +1.24%  "
Set-TextValue "D13" "4.531"
Set-TextValue "E13" "  +0.46%  "
Set-TextValue "D14" "0.5792"
Set-TextValue "E14" "  +0.35%  "
Set-TextValue "D15" "0.000008370"
Set-TextValue "E15" "  -1.50%  "
Set-TextValue "D16" "66.87"
Set-TextValue "E16" "  +3.40%  "
Set-TextValue "D17" "26.328.32"
Set-TextValue "E17" "  -0.01%  "
Set-TextValue "D18" "4.905"
Set-TextValue "E18" "  -0.11%  "
Set-TextValue "D19" "1.008"
Set-TextValue "E19" "  -0.05%  "
Set-TextValue "E20" "  -0.14%  "
Set-TextValue "D21" "193.31"
Set-TextValue "E21" "  +2.00%  "
Set-TextValue "D22" "6.268"
Set-TextValue "E22" "  +1.43%  "
Set-TextValue "D23" "1.008"
Set-TextValue "E23" "  -0.11%  "
Set-TextValue "D24" "149.13"
Set-TextValue "E24" "  +2.93%  "
Set-TextValue "D25" "0.1285"
Set-TextValue "E25" "  +2.03%  "
Set-TextValue "D26" "7.859"
Set-TextValue "E26" "  +1.07%  "
Set-TextValue "D27" "15.82"
Set-TextValue "E27" "  +0.35%  "
Set-TextValue "D28" "1.377"
Set-TextValue "E28" "  +0.99%  "
Set-TextValue "D29" "0.06120"
Set-TextValue "E29" "  -4.69%  "
Set-TextValue "D30" "1.326"
Set-TextValue "E30" "  +0.14%  "
Set-TextValue "D31" "3.605"
Set-TextValue "E31" "  +0.42%  "
Set-TextValue "D32" "3.581"
Set-TextValue "E32" "  -0.04%  "
Set-TextValue "D33" "1.689"
Set-TextValue "E33" "  +1.99%  "
Set-TextValue "D34" "1.033"
Set-TextValue "E34" "  +0.67%  "
Set-TextValue "D35" "0.6188"
Set-TextValue "D36" "2.429"
Set-TextValue "E36" "  +1.01%  "
Set-TextValue "D37" "2.763"
Set-TextValue "E37" "  +0.80%  "
Set-TextValue "D38" "6.236"
Set-TextValue "E38" "  -0.51%  "
Set-TextValue "D39" "0.01639"
Set-TextValue "E39" "  +1.20%  "
Set-TextValue "D40" "1.109.56"
Set-TextValue "E40" "  -0.73%  "
Set-TextValue "D41" "0.8925"
Set-TextValue "E41" "  +2.24%  "
Set-TextValue "D42" "1.013"
Set-TextValue "E42" "  -0.23%  "
Set-TextValue "D43" "100.95"
Set-TextValue "E43" "  +0.50%  "
Set-TextValue "D44" "1.841.30"
Set-TextValue "E44" "  +0.73%  "
Set-TextValue "B45" "Aave"
Set-TextValue "C45" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D45" "57.80"
Set-TextValue "E45" "  +1.62%  "
Set-TextValue "B46" "Frax"
Set-TextValue "C46" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D46" "1.011"
Set-TextValue "E46" "  +0.45%  "
Set-TextValue "B47" "BabyDogeCoin"
Set-TextValue "C47" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D47" "0.00000000105"
Set-TextValue "E47" "  -5.93%  "
Set-TextValue "D48" "8.145"
Set-TextValue "E48" "  -0.24%  "
Set-TextValue "D49" "0.05288"
Set-TextValue "E49" "  +0.43%  "
Set-TextValue "D50" "0.4293"
Set-TextValue "E50" "  -0.10%  "
Set-TextValue "D51" "6.064"
Set-TextValue "E51" "  +0.25%  "
